$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A22").Value = "edit1"
$ws.Range("B22").Value = "riya-morankar"
$ws.Range("C22").Value = "Merged"
$ws.Range("D22").Value = "N/A"

# Assigning a date-like literal ("2025-06-20") directly to .Value would be
# auto-coerced by Excel into a date serial number. To keep it as plain text
# (matching the rest of the "Date" column, which stores these as strings),
# force text format, set the value, then copy the plain style from the row
# above so no stray number-format style is left on the new cell.
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "2025-06-20"
$ws.Range("E21").Copy()
$ws.Range("E22").PasteSpecial(-4122)

$ws.Range("F22").Value = "2afd94510241569eb7fd682e244e8e2f0d248e42"
